$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$data = @(
    @{ Row = 14; Time = "2023-12-06 10:22:13"; Cost = 0.0012 },
    @{ Row = 15; Time = "2023-12-06 10:22:22"; Cost = 0.0002 },
    @{ Row = 16; Time = "2023-12-06 10:22:45"; Cost = 0.0012 },
    @{ Row = 17; Time = "2023-12-06 10:22:55"; Cost = 0.0004 },
    @{ Row = 18; Time = "2023-12-06 10:23:07"; Cost = 0.0004 }
)

foreach ($entry in $data) {
    $r = $entry.Row
    $ws.Cells.Item($r, 1).Value = $entry.Time
    $ws.Cells.Item($r, 2).Value = $entry.Cost
}
